# IST price update 2025-12-23 07:06
#
# A new price snapshot ("2025-12-23 12:32") was captured. It is inserted
# as a new column B, pushing every existing timestamp column (previously
# B:AG) one column to the right (now C:AH). The new column's prices are
# identical to the previous "latest" column (the old column B / now
# column C), since nothing changed price-wise between the two snapshots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B, shifting B:AG -> C:AH.
$ws.Columns("B:B").Insert()

# The inserted column loses its explicit width; restore it to match the
# other price columns (stored width 21 <=> ColumnWidth 20.17 in Excel's
# character-width units).
$ws.Columns("B:B").ColumnWidth = 20.17

# New snapshot timestamp header.
$ws.Range("B1").Value = "2025-12-23 12:32"

# Populate the new column's prices with the latest known price per SKU,
# which is simply what used to be in column B and now lives in column C.
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 3).Value()
}
